# Append three new date rows (5, 6, 7) to column A of Sheet1, matching the
# existing "Date" column's plain-text values (e.g. "2024-09-02" in A2:A4).
#
# Excel's Range.Value setter auto-detects date-shaped strings like
# "2024-10-04" and silently converts them into date serial numbers (with a
# date number format applied). To keep these as literal text - consistent
# with how the rest of column A is stored - format the destination cells as
# Text ("@") before writing the values, then restore the cell style back to
# Normal so no visible/number formatting difference is left behind on the
# new cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A5:A7").NumberFormat = "@"

$ws.Range("A5").Value = "2024-10-04"
$ws.Range("A6").Value = "2024-10-03"
$ws.Range("A7").Value = "2024-10-05"

$ws.Range("A5:A7").Style = "Normal"
